# Updated symbol list — refresh Price (D) and Volume(1h) (E) columns
# for each coin row with the latest scraped values.
#
# The source cells store these figures as literal text (e.g. "312.35",
# "1.34%") rather than numeric/percentage values, so each cell's
# NumberFormat is forced to Text ("@") before the assignment (otherwise
# COM auto-coerces numeric-looking strings into real numbers / applies a
# percent format) and then reset back to the workbook's Normal style so
# no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "312.35"
Set-TextValue $ws.Range("E2") "1.34%"
Set-TextValue $ws.Range("D3") "39.23"
Set-TextValue $ws.Range("E3") "1.15%"
Set-TextValue $ws.Range("D4") "5.132"
Set-TextValue $ws.Range("E4") "0.62%"
Set-TextValue $ws.Range("D5") "0.08129"
Set-TextValue $ws.Range("E5") "-0.03%"
Set-TextValue $ws.Range("D6") "1.985"
Set-TextValue $ws.Range("E6") "0.71%"
Set-TextValue $ws.Range("D7") "4.235"
Set-TextValue $ws.Range("E7") "0.82%"
Set-TextValue $ws.Range("D8") "8.133"
Set-TextValue $ws.Range("E8") "2.50%"
Set-TextValue $ws.Range("D9") "0.9260"
Set-TextValue $ws.Range("E9") "-0.37%"
Set-TextValue $ws.Range("D10") "0.1423"
Set-TextValue $ws.Range("E10") "-1.12%"
Set-TextValue $ws.Range("D11") "0.1928"
Set-TextValue $ws.Range("E11") "-1.65%"
Set-TextValue $ws.Range("D12") "0.09066"
Set-TextValue $ws.Range("E12") "-0.28%"
Set-TextValue $ws.Range("D13") "0.03533"
Set-TextValue $ws.Range("E13") "0.64%"
Set-TextValue $ws.Range("D14") "0.09807"
Set-TextValue $ws.Range("E14") "-0.04%"
Set-TextValue $ws.Range("D15") "0.001390"
Set-TextValue $ws.Range("E15") "-1.23%"
Set-TextValue $ws.Range("D16") "0.006117"
Set-TextValue $ws.Range("E16") "3.23%"
Set-TextValue $ws.Range("D17") "3.762"
Set-TextValue $ws.Range("E17") "4.55%"
Set-TextValue $ws.Range("D18") "3.377"
Set-TextValue $ws.Range("E18") "-0.98%"
Set-TextValue $ws.Range("D19") "0.3454"
Set-TextValue $ws.Range("E19") "0.20%"
Set-TextValue $ws.Range("D20") "0.1312"
Set-TextValue $ws.Range("E20") "-1.74%"
Set-TextValue $ws.Range("D21") "4.638"
Set-TextValue $ws.Range("E21") "-3.84%"
Set-TextValue $ws.Range("D22") "0.2425"
Set-TextValue $ws.Range("E22") "0.90%"
Set-TextValue $ws.Range("D23") "0.04362"
Set-TextValue $ws.Range("E23") "-1.67%"
Set-TextValue $ws.Range("D24") "0.001229"
Set-TextValue $ws.Range("E24") "-0.18%"
Set-TextValue $ws.Range("D25") "0.004801"
Set-TextValue $ws.Range("E25") "-0.97%"
Set-TextValue $ws.Range("D26") "0.0001300"
Set-TextValue $ws.Range("E26") "-0.10%"
Set-TextValue $ws.Range("D27") "0.0004002"
Set-TextValue $ws.Range("E27") "-10.03%"
Set-TextValue $ws.Range("D39") "0.02132"
Set-TextValue $ws.Range("E39") "1.90%"
Set-TextValue $ws.Range("D40") "0.05165"
Set-TextValue $ws.Range("E40") "1.07%"
Set-TextValue $ws.Range("D41") "0.007449"
Set-TextValue $ws.Range("E41") "-0.20%"
Set-TextValue $ws.Range("D42") "0.009773"
Set-TextValue $ws.Range("E42") "-3.09%"
Set-TextValue $ws.Range("D43") "0.1367"
Set-TextValue $ws.Range("E43") "0.29%"
Set-TextValue $ws.Range("D44") "0.002131"
Set-TextValue $ws.Range("E44") "-0.57%"
Set-TextValue $ws.Range("D45") "0.009757"
Set-TextValue $ws.Range("E45") "-6.85%"
Set-TextValue $ws.Range("D46") "0.00006390"
Set-TextValue $ws.Range("E46") "2.45%"
Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "-0.10%"
Set-TextValue $ws.Range("D48") "0.001000"
Set-TextValue $ws.Range("E48") "-37.56%"
Set-TextValue $ws.Range("D49") "0.002555"
Set-TextValue $ws.Range("E49") "-16.72%"
Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("E50") "-0.10%"
Set-TextValue $ws.Range("D51") "0.0002001"
Set-TextValue $ws.Range("E51") "-0.10%"
